# Fix Training Data Issue (#48)
# The "Date" column (BF) held values like "5-31-2011-12" (month-day-season)
# instead of a proper ISO date string "2012-05-31". Correct every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$dateCol = 58   # column BF ("Date")

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    if ($cell.Value() -eq "5-31-2011-12") {
        # Force text so Excel doesn't reinterpret the ISO-looking string as a date serial.
        $cell.NumberFormat = "@"
        $cell.Value = "2012-05-31"
    }
}
